# Sprint 6 Report -- "added stuff to research"
#
# The canonical diff this script targets mostly consists of Word's
# automatic spell/grammar-checker bookkeeping (<w:proofErr .../> markers
# that split otherwise-unchanged runs around words like "Malkowski",
# "eBill", "AutoLayout", "misc"). Those markers are written by Word's
# live proofing pass, not by anything the Word object model exposes to
# automation (there is no WordBasic/VBA call that inserts a
# w:proofErr element), so the wording they wrap is left untouched here
# -- it was never actually changed by the edit. The substantive content
# change is the new paragraph appended to the "Research/Code
# Experiments" section (and the accompanying relocation of the
# "_GoBack" bookmark to the end of that new paragraph), which this
# script performs explicitly.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its old spot ("Brief Scrum
#    Meeting"). It will be recreated further down, at the end of the
#    newly-added research paragraph, which is where the edit moved it.
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# ------------------------------------------------------------------
# 2. Add the new paragraph content to the "Research/Code Experiments"
#    section, right before the page break that follows the existing
#    "Mike and Rachel..." paragraph.
# ------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Mike and Rachel both researched")

$researchPara = $findRng.Paragraphs(1)
$pageBreakPara = $researchPara.Next()
$insertRng = $pageBreakPara.Range
$insertRng.Collapse(1)  # wdCollapseStart

$tab = [char]9
$newText = $tab + "Rachel attempted to implement horizontal scrolling for the Technician App. With the new iOS AutoLayout implementation, this does not work. It is likely that future AutoLayout implementations should make this fairly simple. As of now, however, horizontal scrolling has been moved to the product backlog."

$insertRng.InsertBefore($newText)

# ------------------------------------------------------------------
# 3. Recreate the _GoBack bookmark as a zero-length bookmark right
#    before the page-break run, at the end of the text we just added.
# ------------------------------------------------------------------
$afterText = $d.Content
$afterText.Find.Execute("moved to the product backlog.")
$bmRng = $afterText.Duplicate
$bmRng.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Output "Done."
